# Add two additional applications (Happiness Analysis, Weather Forecast Web
# App) as new rows 11 and 12 in the "Python Programs" sheet, pushing the
# existing placeholder row down to row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Python Programs")

# ---------------------------------------------------------------------
# 1. Make room: insert one new blank row at position 12. This shifts the
#    old row 12 (placeholder) down to row 13, while old row 11
#    (placeholder, about to be overwritten) remains at row 11 for now.
# ---------------------------------------------------------------------
$ws.Rows("12:12").Insert()

# ---------------------------------------------------------------------
# 2. Give the new row 12 the same cell formatting pattern used by the
#    other "application" rows (A/B/D = normal style, C = hyperlink
#    style with wrap) by copying formats only from row 10, so no new
#    style entries are created.
# ---------------------------------------------------------------------
$ws.Range("A10:D10").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Overwrite row 11 in place with the "Happiness Analysis" entry.
#    (Its existing formatting already matches A=4,B=3,D=2; only column C
#    needs the hyperlink-with-wrap style, copied from an existing C
#    cell that already uses it so the style table stays unchanged.)
# ---------------------------------------------------------------------
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Happiness Analysis"
$ws.Range("B11").Value = "Several factors, both economic and social, increase the happiness of populations in countries. The study data gathered for several contries is part of this application. The streamlit application allows the user to plot graphs of these factors and overall happiness."
$ws.Range("C11").Value = "https://github.com/valenpendragon/happiness-web-app"
$ws.Range("D11").Value = "5.png"
$ws.Rows("11:11").RowHeight = 60

# ---------------------------------------------------------------------
# 4. Fill the new row 12 with the "Weather Forecast Web App" entry.
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "Weather Forecast Web App"
$ws.Range("B12").Value = "This streamlit appliction gathers up to 5 days of forecast sky conditions and temperature data for a city they wish to view."
$ws.Range("C12").Value = "https://github.com/valenpendragon/weather-forecast-web-app"
$ws.Range("D12").Value = "5.png"
$ws.Rows("12:12").RowHeight = 30

# ---------------------------------------------------------------------
# 5. Add the hyperlinks for the two new URL cells.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C11"), "https://github.com/valenpendragon/happiness-web-app")
$ws.Hyperlinks.Add($ws.Range("C12"), "https://github.com/valenpendragon/weather-forecast-web-app")

# ---------------------------------------------------------------------
# 6. Restore the view so the newly added row is visible/selected, as it
#    would be right after typing the data in.
# ---------------------------------------------------------------------
$ws.Range("C11").Select()
